$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 69
$ws.Range("J9").Value = 70
$ws.Range("L9").Value = 70
$ws.Range("N9").Value = -408
$ws.Range("H29").Value = 234.5
$ws.Range("I29").Value = 234.5
$ws.Range("K29").Value = 703.5
$ws.Range("M29").Value = -422.5
$ws.Range("H38").Value = 2332.3
$ws.Range("I38").Value = 210.8
$ws.Range("J38").Value = 4453.8
$ws.Range("K38").Value = 632.4000000000001
$ws.Range("L38").Value = 13361.4
$ws.Range("M38").Value = -260.4000000000001
$ws.Range("N38").Value = -14105.4
$ws.Range("H58").Value = 508.57144
$ws.Range("I58").Value = 176.66667
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 530.00001
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -380.00001
$ws.Range("N58").Value = -7800
$ws.Range("H70").Value = 1722.2222
$ws.Range("J70").Value = 2750
$ws.Range("L70").Value = 8250
$ws.Range("N70").Value = -8790
$ws.Range("H73").Value = 1722.2222
$ws.Range("J73").Value = 2750
$ws.Range("L73").Value = 8250
$ws.Range("N73").Value = -10122
$ws.Range("H92").Value = 251.42105
$ws.Range("I92").Value = 228.57143
$ws.Range("K92").Value = 228.57143
$ws.Range("M92").Value = 1019.42857
$ws.Range("H98").Value = 995.44446
$ws.Range("I98").Value = 995.44446
$ws.Range("K98").Value = 995.44446
$ws.Range("M98").Value = 502.55554
$ws.Range("H113").Value = 43495188
$ws.Range("I113").Value = 52635600
$ws.Range("K113").Value = 52635600
$ws.Range("M113").Value = -52632346
$ws.Range("H122").Value = 995.44446
$ws.Range("I122").Value = 995.44446
$ws.Range("K122").Value = 2986.33338
$ws.Range("M122").Value = -536.33338
$ws.Range("H138").Value = 2442
$ws.Range("J138").Value = 3311.074
$ws.Range("L138").Value = 9933.222
$ws.Range("N138").Value = -20213.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 205950
$ws.Range("I45").Value = 337416.66
$ws.Range("J45").Value = 8750
$ws.Range("K45").Value = 337416.66
$ws.Range("L45").Value = 8750
$ws.Range("M45").Value = -337039.66
$ws.Range("N45").Value = -9504
$ws.Range("H61").Value = 10502.842
$ws.Range("I61").Value = 10502.842
$ws.Range("K61").Value = 10502.842
$ws.Range("M61").Value = -10290.842
$ws.Range("H122").Value = 1147.8889
$ws.Range("I122").Value = 1150.7646
$ws.Range("K122").Value = 3452.2938
$ws.Range("M122").Value = -1002.2938
$ws.Range("H136").Value = 10502.842
$ws.Range("I136").Value = 10502.842
$ws.Range("K136").Value = 31508.526
$ws.Range("M136").Value = -28958.526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10372.5
$ws.Range("I86").Value = 9995
$ws.Range("K86").Value = 9995
$ws.Range("M86").Value = -8872
$ws.Range("H89").Value = 10372.5
$ws.Range("I89").Value = 9995
$ws.Range("K89").Value = 49975
$ws.Range("M89").Value = -44359
$ws.Range("H94").Value = 1684.875
$ws.Range("J94").Value = 1615.6666
$ws.Range("L94").Value = 1615.6666
$ws.Range("N94").Value = -2517.6666
$ws.Range("H97").Value = 24767.25
$ws.Range("I97").Value = 16000
$ws.Range("J97").Value = 26019.715
$ws.Range("K97").Value = 16000
$ws.Range("L97").Value = 26019.715
$ws.Range("M97").Value = -15009
$ws.Range("N97").Value = -28001.715
$ws.Range("H99").Value = 4918.8667
$ws.Range("J99").Value = 5074.5
$ws.Range("L99").Value = 5074.5
$ws.Range("N99").Value = -8070.5
$ws.Range("H107").Value = 310.33334
$ws.Range("J107").Value = 184.75
$ws.Range("L107").Value = 184.75
$ws.Range("N107").Value = -4024.75
$ws.Range("H122").Value = 2872.2173
$ws.Range("I122").Value = 3073.1
$ws.Range("K122").Value = 9219.299999999999
$ws.Range("M122").Value = -6769.299999999999
$ws.Range("H126").Value = 4918.8667
$ws.Range("J126").Value = 5074.5
$ws.Range("L126").Value = 15223.5
$ws.Range("N126").Value = -20163.5
$ws.Range("H141").Value = 124296.92
$ws.Range("I141").Value = 60997
$ws.Range("K141").Value = 60997
$ws.Range("M141").Value = -55817

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 849.5
$ws.Range("I3").Value = 849.5
$ws.Range("K3").Value = 2548.5
$ws.Range("M3").Value = -2436.5
$ws.Range("H113").Value = 1670.6
$ws.Range("I113").Value = 585.25
$ws.Range("K113").Value = 1755.75
$ws.Range("M113").Value = 414.25
$ws.Range("H137").Value = 11704.75
$ws.Range("I137").Value = 7439.6665
$ws.Range("K137").Value = 22318.9995
$ws.Range("M137").Value = -17218.9995
$ws.Range("H140").Value = 2125
$ws.Range("I140").Value = 2125
$ws.Range("K140").Value = 6375
$ws.Range("M140").Value = -1195

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1747.5
$ws.Range("I122").Value = 1722.4375
$ws.Range("K122").Value = 5167.3125
$ws.Range("M122").Value = -2717.3125
$ws.Range("H127").Value = 55554
$ws.Range("J127").Value = 55554
$ws.Range("L127").Value = 55554
$ws.Range("N127").Value = -65474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6336.35
$ws.Range("I46").Value = 6732.1113
$ws.Range("K46").Value = 6732.1113
$ws.Range("M46").Value = -6544.1113
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7195
$ws.Range("J15").Value = 7195
$ws.Range("L15").Value = 7195
$ws.Range("N15").Value = -7771
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9490
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 35052.5
$ws.Range("I70").Value = 30000
$ws.Range("J70").Value = 40105
$ws.Range("K70").Value = 30000
$ws.Range("L70").Value = 40105
$ws.Range("M70").Value = -29685
$ws.Range("N70").Value = -40735
$ws.Range("H73").Value = 35052.5
$ws.Range("I73").Value = 30000
$ws.Range("J73").Value = 40105
$ws.Range("K73").Value = 30000
$ws.Range("L73").Value = 40105
$ws.Range("M73").Value = -28908
$ws.Range("N73").Value = -42289
$ws.Range("H107").Value = 1207.4667
$ws.Range("I107").Value = 575.375
$ws.Range("J107").Value = 1929.8572
$ws.Range("K107").Value = 1726.125
$ws.Range("L107").Value = 5789.571599999999
$ws.Range("M107").Value = 193.875
$ws.Range("N107").Value = -9629.571599999999
$ws.Range("H122").Value = 3505.8572
$ws.Range("I122").Value = 2500.077
$ws.Range("K122").Value = 7500.231000000001
$ws.Range("M122").Value = -5050.231000000001
$ws.Range("H132").Value = 4559.278
$ws.Range("I132").Value = 3829.6333
$ws.Range("K132").Value = 11488.8999
$ws.Range("M132").Value = -8958.8999
